# Auto-generated script applying market-price / profit recalculation updates
# to the Sheets workbook (ALC, ARM, BSM, CRP, GSM, LTW, WVR tables).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6286.4
$ws.Range("I40").Value = 3816
$ws.Range("J40").Value = 8307.637000000001
$ws.Range("K40").Value = 3816
$ws.Range("L40").Value = 8307.637000000001
$ws.Range("M40").Value = -3641
$ws.Range("N40").Value = -8657.637000000001
$ws.Range("H64").Value = 6999.2812
$ws.Range("I64").Value = 5682.625
$ws.Range("J64").Value = 8315.9375
$ws.Range("K64").Value = 5682.625
$ws.Range("L64").Value = 8315.9375
$ws.Range("M64").Value = -5434.625
$ws.Range("N64").Value = -8811.9375
$ws.Range("H67").Value = 6999.2812
$ws.Range("I67").Value = 5682.625
$ws.Range("J67").Value = 8315.9375
$ws.Range("K67").Value = 5682.625
$ws.Range("L67").Value = 8315.9375
$ws.Range("M67").Value = -4824.625
$ws.Range("N67").Value = -10031.9375
$ws.Range("H74").Value = 7224.25
$ws.Range("I74").Value = 4790.4165
$ws.Range("K74").Value = 4790.4165
$ws.Range("M74").Value = -3854.4165
$ws.Range("H77").Value = 7224.25
$ws.Range("I77").Value = 4790.4165
$ws.Range("K77").Value = 23952.0825
$ws.Range("M77").Value = -19272.0825
$ws.Range("H80").Value = 2822.2068
$ws.Range("I80").Value = 1536.6666
$ws.Range("J80").Value = 4199.5713
$ws.Range("K80").Value = 4609.9998
$ws.Range("L80").Value = 12598.7139
$ws.Range("M80").Value = -3611.9998
$ws.Range("N80").Value = -14594.7139
$ws.Range("H83").Value = 2822.2068
$ws.Range("I83").Value = 1536.6666
$ws.Range("J83").Value = 4199.5713
$ws.Range("K83").Value = 13829.9994
$ws.Range("L83").Value = 37796.14169999999
$ws.Range("M83").Value = -8837.999400000001
$ws.Range("N83").Value = -47780.14169999999
$ws.Range("H92").Value = 4808758.5
$ws.Range("I92").Value = 785.7778
$ws.Range("K92").Value = 785.7778
$ws.Range("M92").Value = 462.2222
$ws.Range("H131").Value = 6421.636
$ws.Range("I131").Value = 4376.5
$ws.Range("J131").Value = 15624.75
$ws.Range("K131").Value = 13129.5
$ws.Range("L131").Value = 46874.25
$ws.Range("M131").Value = -8089.5
$ws.Range("N131").Value = -56954.25
$ws.Range("H137").Value = 2812.0833
$ws.Range("I137").Value = 3029.4443
$ws.Range("K137").Value = 9088.332900000001
$ws.Range("M137").Value = -6538.332900000001
$ws.Range("H138").Value = 2847.4314
$ws.Range("J138").Value = 3012.7954
$ws.Range("L138").Value = 9038.386200000001
$ws.Range("N138").Value = -19318.3862
$ws.Range("H141").Value = 4666.4346
$ws.Range("I141").Value = 4635.619
$ws.Range("K141").Value = 13906.857
$ws.Range("M141").Value = -8726.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1275
$ws.Range("I3").Value = 1275
$ws.Range("K3").Value = 1275
$ws.Range("M3").Value = -1160
$ws.Range("H63").Value = 2239.9375
$ws.Range("J63").Value = 1806.3334
$ws.Range("L63").Value = 1806.3334
$ws.Range("N63").Value = -3178.3334
$ws.Range("H66").Value = 2239.9375
$ws.Range("J66").Value = 1806.3334
$ws.Range("L66").Value = 9031.666999999999
$ws.Range("N66").Value = -15895.667
$ws.Range("H110").Value = 1767.5264
$ws.Range("I110").Value = 1911.7646
$ws.Range("J110").Value = 541.5
$ws.Range("K110").Value = 1911.7646
$ws.Range("L110").Value = 541.5
$ws.Range("M110").Value = 133.2354
$ws.Range("N110").Value = -4631.5
$ws.Range("H132").Value = 4444.65
$ws.Range("I132").Value = 3565.2856
$ws.Range("K132").Value = 10695.8568
$ws.Range("M132").Value = -8165.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10675.875
$ws.Range("I86").Value = 11706
$ws.Range("J86").Value = 10528.714
$ws.Range("K86").Value = 11706
$ws.Range("L86").Value = 10528.714
$ws.Range("M86").Value = -10583
$ws.Range("N86").Value = -12774.714
$ws.Range("H89").Value = 10675.875
$ws.Range("I89").Value = 11706
$ws.Range("J89").Value = 10528.714
$ws.Range("K89").Value = 58530
$ws.Range("L89").Value = 52643.57
$ws.Range("M89").Value = -52914
$ws.Range("N89").Value = -63875.57
$ws.Range("H134").Value = 2642.0732
$ws.Range("I134").Value = 1606.5938
$ws.Range("J134").Value = 6323.778
$ws.Range("K134").Value = 4819.7814
$ws.Range("L134").Value = 18971.334
$ws.Range("M134").Value = -2284.7814
$ws.Range("N134").Value = -24041.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2671.08
$ws.Range("I31").Value = 2024.5116
$ws.Range("K31").Value = 2024.5116
$ws.Range("M31").Value = -1729.5116
$ws.Range("H34").Value = 2671.08
$ws.Range("I34").Value = 2024.5116
$ws.Range("K34").Value = 2024.5116
$ws.Range("M34").Value = -1822.5116
$ws.Range("H62").Value = 4554.5386
$ws.Range("I62").Value = 5399.625
$ws.Range("J62").Value = 3202.4
$ws.Range("K62").Value = 5399.625
$ws.Range("L62").Value = 3202.4
$ws.Range("M62").Value = -4775.625
$ws.Range("N62").Value = -4450.4
$ws.Range("H65").Value = 4554.5386
$ws.Range("I65").Value = 5399.625
$ws.Range("J65").Value = 3202.4
$ws.Range("K65").Value = 26998.125
$ws.Range("L65").Value = 16012
$ws.Range("M65").Value = -23878.125
$ws.Range("N65").Value = -22252
$ws.Range("H99").Value = 9761771
$ws.Range("I99").Value = 2038562.1
$ws.Range("K99").Value = 2038562.1
$ws.Range("M99").Value = -2037064.1
$ws.Range("H105").Value = 596.925
$ws.Range("I105").Value = 457.2414
$ws.Range("J105").Value = 965.1818
$ws.Range("K105").Value = 457.2414
$ws.Range("L105").Value = 965.1818
$ws.Range("M105").Value = 1289.7586
$ws.Range("N105").Value = -4459.1818
$ws.Range("H126").Value = 9761771
$ws.Range("I126").Value = 2038562.1
$ws.Range("K126").Value = 6115686.300000001
$ws.Range("M126").Value = -6113216.300000001
$ws.Range("H132").Value = 5179.5
$ws.Range("J132").Value = 6444
$ws.Range("L132").Value = 19332
$ws.Range("N132").Value = -24392
$ws.Range("H134").Value = 4486.087
$ws.Range("I134").Value = 2655.9285
$ws.Range("K134").Value = 7967.7855
$ws.Range("M134").Value = -5432.7855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 296.42856
$ws.Range("I13").Value = 282.83334
$ws.Range("J13").Value = 306.625
$ws.Range("K13").Value = 282.83334
$ws.Range("L13").Value = 306.625
$ws.Range("M13").Value = -143.83334
$ws.Range("N13").Value = -584.625
$ws.Range("H97").Value = 12500410
$ws.Range("I97").Value = 396.7857
$ws.Range("J97").Value = 41667108
$ws.Range("K97").Value = 396.7857
$ws.Range("L97").Value = 41667108
$ws.Range("M97").Value = 99.21429999999998
$ws.Range("N97").Value = -41668100
$ws.Range("H132").Value = 4480.433
$ws.Range("I132").Value = 4127.8945
$ws.Range("J132").Value = 5089.364
$ws.Range("K132").Value = 12383.6835
$ws.Range("L132").Value = 15268.092
$ws.Range("M132").Value = -9853.683500000001
$ws.Range("N132").Value = -20328.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 37999.75
$ws.Range("J6").Value = 37999.75
$ws.Range("L6").Value = 37999.75
$ws.Range("N6").Value = -38223.75
$ws.Range("H68").Value = 269389.3
$ws.Range("I68").Value = 255000
$ws.Range("K68").Value = 255000
$ws.Range("M68").Value = -254251
$ws.Range("H71").Value = 269389.3
$ws.Range("I71").Value = 255000
$ws.Range("K71").Value = 1275000
$ws.Range("M71").Value = -1271256

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 721.53845
$ws.Range("I107").Value = 508
$ws.Range("K107").Value = 1524
$ws.Range("M107").Value = 396
$ws.Range("H113").Value = 1024.762
$ws.Range("I113").Value = 708.06665
$ws.Range("J113").Value = 1816.5
$ws.Range("K113").Value = 2124.19995
$ws.Range("L113").Value = 5449.5
$ws.Range("M113").Value = 45.80004999999983
$ws.Range("N113").Value = -9789.5
$ws.Range("H126").Value = 1547.5555
$ws.Range("I126").Value = 1515.75
$ws.Range("K126").Value = 4547.25
$ws.Range("M126").Value = -2077.25

Write-Host "Applied 203 cell updates across 7 sheets"
